{"js": "// Replace the date line and each three-digit-divided-by-one-digit problem text\n// with its updated value, one unique search-and-replace per run.\nconst replacements = [\n  [\"2024-03-04 Monday\", \"2024-03-05 Tuesday\"],\n  [\"300\u00f73=100, 0\", \"474\u00f74=118, 2\"],\n  [\"517\u00f78=64, 5\", \"429\u00f72=214, 1\"],\n  [\"133\u00f73=44, 1\", \"628\u00f77=89, 5\"],\n  [\"199\u00f73=66, 1\", \"580\u00f75=116, 0\"],\n  [\"858\u00f76=143, 0\", \"220\u00f77=31, 3\"],\n  [\"770\u00f76=128, 2\", \"301\u00f72=150, 1\"],\n  [\"392\u00f73=130, 2\", \"648\u00f77=92, 4\"],\n  [\"805\u00f77=115, 0\", \"947\u00f78=118, 3\"],\n  [\"141\u00f78=17, 5\", \"360\u00f79=40, 0\"],\n  [\"682\u00f72=341, 0\", \"411\u00f79=45, 6\"],\n  [\"820\u00f74=205, 0\", \"717\u00f79=79, 6\"],\n  [\"984\u00f76=164, 0\", \"416\u00f72=208, 0\"],\n  [\"410\u00f75=82, 0\", \"918\u00f78=114, 6\"],\n  [\"900\u00f72=450, 0\", \"609\u00f78=76, 1\"],\n  [\"608\u00f75=121, 3\", \"452\u00f73=150, 2\"],\n  [\"314\u00f76=52, 2\", \"860\u00f72=430, 0\"],\n  [\"524\u00f72=262, 0\", \"665\u00f76=110, 5\"],\n  [\"940\u00f74=235, 0\", \"786\u00f78=98, 2\"],\n  [\"350\u00f72=175, 0\", \"335\u00f79=37, 2\"],\n  [\"329\u00f75=65, 4\", \"447\u00f78=55, 7\"],\n  [\"954\u00f76=159, 0\", \"539\u00f74=134, 3\"],\n  [\"256\u00f78=32, 0\", \"189\u00f76=31, 3\"],\n  [\"606\u00f76=101, 0\", \"464\u00f77=66, 2\"],\n  [\"224\u00f76=37, 2\", \"381\u00f77=54, 3\"],\n  [\"771\u00f76=128, 3\", \"495\u00f77=70, 5\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each three-digit-divided-by-one-digit problem text\n# with its updated value, using Word's Find/Replace on the document content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-03-04 Monday', '2024-03-05 Tuesday'),\n    @('300\u00f73=100, 0', '474\u00f74=118, 2'),\n    @('517\u00f78=64, 5', '429\u00f72=214, 1'),\n    @('133\u00f73=44, 1', '628\u00f77=89, 5'),\n    @('199\u00f73=66, 1', '580\u00f75=116, 0'),\n    @('858\u00f76=143, 0', '220\u00f77=31, 3'),\n    @('770\u00f76=128, 2', '301\u00f72=150, 1'),\n    @('392\u00f73=130, 2', '648\u00f77=92, 4'),\n    @('805\u00f77=115, 0', '947\u00f78=118, 3'),\n    @('141\u00f78=17, 5', '360\u00f79=40, 0'),\n    @('682\u00f72=341, 0', '411\u00f79=45, 6'),\n    @('820\u00f74=205, 0', '717\u00f79=79, 6'),\n    @('984\u00f76=164, 0', '416\u00f72=208, 0'),\n    @('410\u00f75=82, 0', '918\u00f78=114, 6'),\n    @('900\u00f72=450, 0', '609\u00f78=76, 1'),\n    @('608\u00f75=121, 3', '452\u00f73=150, 2'),\n    @('314\u00f76=52, 2', '860\u00f72=430, 0'),\n    @('524\u00f72=262, 0', '665\u00f76=110, 5'),\n    @('940\u00f74=235, 0', '786\u00f78=98, 2'),\n    @('350\u00f72=175, 0', '335\u00f79=37, 2'),\n    @('329\u00f75=65, 4', '447\u00f78=55, 7'),\n    @('954\u00f76=159, 0', '539\u00f74=134, 3'),\n    @('256\u00f78=32, 0', '189\u00f76=31, 3'),\n    @('606\u00f76=101, 0', '464\u00f77=66, 2'),\n    @('224\u00f76=37, 2', '381\u00f77=54, 3'),\n    @('771\u00f76=128, 3', '495\u00f77=70, 5'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
